$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift word/translation pairs: rows 320-566 now hold the word that used to be
# one row below (the old row-343 entry "insinoori/engineer" was removed from
# the middle of the list), and two new entries are appended at the end of the
# range (rows 565-566), each scored correct (1).
$updates = @(
    @{Row=320; A='ihminen'; B='human'}
    @{Row=321; A='hyvin'; B='well (adverb)'}
    @{Row=322; A='sohva'; B='sofa'}
    @{Row=323; A='kanadalainen'; B='Canadian (adjective)'}
    @{Row=324; A='puraista'; B='to bite'}
    @{Row=325; A='omena'; B='apple'}
    @{Row=326; A='korvapuusti'; B='cinnamon roll'}
    @{Row=327; A='edes'; B='even'}
    @{Row=328; A='ehkä'; B='maybe'}
    @{Row=329; A='kasvi'; B='plant'}
    @{Row=330; A='onnellinen'; B='happy'}
    @{Row=331; A='litra'; B='liter'}
    @{Row=332; A='Pariisi'; B='Paris'}
    @{Row=333; A='kilo'; B='kilogram'}
    @{Row=334; A='poikaystävä'; B='boyfriend'}
    @{Row=335; A='suomalainen'; B='Finnish person'}
    @{Row=336; A='onni'; B='happiness'}
    @{Row=337; A='nämä'; B='these'}
    @{Row=338; A='vaari'; B='grandpa'}
    @{Row=339; A='karjalanpiirakka'; B='Karelian pasty'}
    @{Row=340; A='saksa'; B='German language'}
    @{Row=341; A='myös'; B='also'}
    @{Row=342; A='kyllä'; B='yes'}
    @{Row=343; A='kana'; B='chicken'}
    @{Row=344; A='laulaja'; B='singer'}
    @{Row=345; A='sokeri'; B='sugar'}
    @{Row=346; A='viisi'; B='five'}
    @{Row=347; A='tajuta'; B='to realize'}
    @{Row=348; A='tähti'; B='star'}
    @{Row=349; A='ahma'; B='wolwerine'}
    @{Row=350; A='valkoinen'; B='white'}
    @{Row=351; A='karhu'; B='bear'}
    @{Row=352; A='kellari'; B='basement'}
    @{Row=353; A='Elsa'; B='Elsa'}
    @{Row=354; A='he'; B='they'}
    @{Row=355; A='susi'; B='wolf'}
    @{Row=356; A='kaupunki'; B='city'}
    @{Row=357; A='morsian'; B='bride'}
    @{Row=358; A='Australia'; B='Australia'}
    @{Row=359; A='afrikkalainen'; B='African'}
    @{Row=360; A='peli'; B='game'}
    @{Row=361; A='tosi'; B='really (adverb of degree, spoken language)'}
    @{Row=362; A='käärme'; B='snake'}
    @{Row=363; A='kirjoittaa'; B='to write'}
    @{Row=364; A='koira'; B='dog'}
    @{Row=365; A='Pyry'; B='Pyry'}
    @{Row=366; A='mukava'; B='comfortable'}
    @{Row=367; A='Suomi'; B='Finland'}
    @{Row=368; A='ruotsalainen'; B='Swedish person'}
    @{Row=369; A='herkullinen'; B='delicious'}
    @{Row=370; A='velho'; B='wizard'}
    @{Row=371; A='ovat'; B='(they) are'}
    @{Row=372; A='grilli'; B='grill'}
    @{Row=373; A='pyöreä'; B='round'}
    @{Row=374; A='kolme'; B='three'}
    @{Row=375; A='tuo'; B='that (adjective)'}
    @{Row=376; A='hän'; B='s/he'}
    @{Row=377; A='itkeä'; B='to cry'}
    @{Row=378; A='tanskalainen'; B='Danish person'}
    @{Row=379; A='makkara'; B='sausage'}
    @{Row=380; A='jotain'; B='something'}
    @{Row=381; A='sinivalkoinen'; B='blue and white'}
    @{Row=382; A='loppu'; B='out of'}
    @{Row=383; A='espanja'; B='Spanish language'}
    @{Row=384; A='kuusi'; B='spruce'}
    @{Row=385; A='väärin'; B='incorrect (adverb)'}
    @{Row=386; A='jo'; B='already'}
    @{Row=387; A='aika'; B='time'}
    @{Row=388; A='mämmi'; B='mämmi'}
    @{Row=389; A='opettaja'; B='teacher'}
    @{Row=390; A='tyttöystävä'; B='girlfriend'}
    @{Row=391; A='hevonen'; B='horse'}
    @{Row=392; A='Oslo'; B='Oslo'}
    @{Row=393; A='oikea'; B='right (adjective)'}
    @{Row=394; A='-kä'; B='nor'}
    @{Row=395; A='kaukana'; B='far away'}
    @{Row=396; A='missä'; B='where'}
    @{Row=397; A='vielä'; B='still'}
    @{Row=398; A='akku'; B='battery'}
    @{Row=399; A='vadelma'; B='raspberry'}
    @{Row=400; A='puoli'; B='half'}
    @{Row=401; A='hindi'; B='Hindi (language)'}
    @{Row=402; A='ketkä'; B='who (are the people)'}
    @{Row=403; A='pehmeä'; B='soft'}
    @{Row=404; A='voi voi'; B='oh dear (~butter)'}
    @{Row=405; A='kaikki'; B='everyone'}
    @{Row=406; A='laturi'; B='charger'}
    @{Row=407; A='kuin'; B='as (comparison)'}
    @{Row=408; A='Skandinavia'; B='Skandinavia'}
    @{Row=409; A='sauna'; B='sauna'}
    @{Row=410; A='sisu'; B='true grit'}
    @{Row=411; A='kanadalainen'; B='Canadian person'}
    @{Row=412; A='radio'; B='radio'}
    @{Row=413; A='ihahaa'; B='neigh (horse sound)'}
    @{Row=414; A='netti'; B='net (IT)'}
    @{Row=415; A='murre'; B='dialect'}
    @{Row=416; A='tai'; B='or (inclusive)'}
    @{Row=417; A='kirkas'; B='bright'}
    @{Row=418; A='likainen'; B='dirty'}
    @{Row=419; A='Kanada'; B='Canada'}
    @{Row=420; A='vaikea'; B='difficult'}
    @{Row=421; A='vaarallinen'; B='dangerous'}
    @{Row=422; A='viisas'; B='wise'}
    @{Row=423; A='kameli'; B='camel'}
    @{Row=424; A='terve'; B='hello'}
    @{Row=425; A='kiinni'; B='closed'}
    @{Row=426; A='miau'; B='meow (cat sound)'}
    @{Row=427; A='Viro'; B='Estonia'}
    @{Row=428; A='naimisissa'; B='married'}
    @{Row=429; A='oopperatalo'; B='opera house'}
    @{Row=430; A='onnea'; B='good luck'}
    @{Row=431; A='kenguru'; B='kangaroo'}
    @{Row=432; A='monumentti'; B='monument'}
    @{Row=433; A='sininen'; B='blue'}
    @{Row=434; A='voida'; B='may, can, to be able to'}
    @{Row=435; A='röh'; B='oink'}
    @{Row=436; A='nyt'; B='now'}
    @{Row=437; A='kiitos'; B='thank you'}
    @{Row=438; A='au'; B='ouch'}
    @{Row=439; A='myydä'; B='to sell'}
    @{Row=440; A='paljonko'; B='how much'}
    @{Row=441; A='veri'; B='blood'}
    @{Row=442; A='mauste'; B='spice'}
    @{Row=443; A='oletko'; B='are (you singular in question)'}
    @{Row=444; A='sima'; B='mead'}
    @{Row=445; A='islantilainen'; B='Icelandic person'}
    @{Row=446; A='glögi'; B='glögi'}
    @{Row=447; A='saamelainen'; B='Sámi person'}
    @{Row=448; A='venäläinen'; B='Russian person'}
    @{Row=449; A='lemmikki'; B='pet'}
    @{Row=450; A='osata'; B='to know how to'}
    @{Row=451; A='kuinka'; B='how'}
    @{Row=452; A='limonadi'; B='soda pop'}
    @{Row=453; A='koko'; B='the whole'}
    @{Row=454; A='murista'; B='to growl'}
    @{Row=455; A='lämmin'; B='warm'}
    @{Row=456; A='olla ikävä'; B='to miss'}
    @{Row=457; A='Berliini'; B='Berlin'}
    @{Row=458; A='sulaa'; B='to melt'}
    @{Row=459; A='ruoka'; B='food'}
    @{Row=460; A='jolla'; B='who/which (has) (relative)'}
    @{Row=461; A='auki'; B='open'}
    @{Row=462; A='nopeasti'; B='quickly'}
    @{Row=463; A='Islanti'; B='Iceland'}
    @{Row=464; A='nauraa'; B='to laugh'}
    @{Row=465; A='Tyyne'; B='Tyyne'}
    @{Row=466; A='rahka'; B='quark'}
    @{Row=467; A='stadion'; B='stadium'}
    @{Row=468; A='kahdeksan'; B='eight'}
    @{Row=469; A='pupu'; B='bunny'}
    @{Row=470; A='tuhma'; B='naughty'}
    @{Row=471; A='marja'; B='berry'}
    @{Row=472; A='hauska'; B='funny'}
    @{Row=473; A='sydän'; B='heart'}
    @{Row=474; A='taas'; B='again'}
    @{Row=475; A='meri'; B='sea'}
    @{Row=476; A='lehti'; B='leaf'}
    @{Row=477; A='auto'; B='car'}
    @{Row=478; A='varma'; B='sure, certain'}
    @{Row=479; A='sama'; B='same'}
    @{Row=480; A='minuutti'; B='minute'}
    @{Row=481; A='korea'; B='Korean language'}
    @{Row=482; A='kivi'; B='stone, rock'}
    @{Row=483; A='kai'; B='maybe, I guess'}
    @{Row=484; A='yrittää'; B='to try'}
    @{Row=485; A='älykäs'; B='intelligent'}
    @{Row=486; A='metsä'; B='forest'}
    @{Row=487; A='poika'; B='boy'}
    @{Row=488; A='tabletti'; B='tablet'}
    @{Row=489; A='juosta'; B='to run'}
    @{Row=490; A='mummo'; B='grandma'}
    @{Row=491; A='suklaa'; B='chocolate'}
    @{Row=492; A='hei'; B='hi (greeting)'}
    @{Row=493; A='kiisseli'; B='kissel'}
    @{Row=494; A='syötävä'; B='edible'}
    @{Row=495; A='uusi'; B='new'}
    @{Row=496; A='karkki'; B='candy'}
    @{Row=497; A='kantele'; B='kantele'}
    @{Row=498; A='ranska'; B='French language'}
    @{Row=499; A='haarukka'; B='fork'}
    @{Row=500; A='vihreä'; B='green'}
    @{Row=501; A='komea'; B='handsome'}
    @{Row=502; A='viettää'; B='to spend (time, vacation)'}
    @{Row=503; A='monta'; B='many (partitive)'}
    @{Row=504; A='kiltti'; B='well-behaved'}
    @{Row=505; A='tasan'; B='even, exactly (time)'}
    @{Row=506; A='kuppi'; B='cup (of)'}
    @{Row=507; A='kartta'; B='map'}
    @{Row=508; A='Lontoo'; B='London'}
    @{Row=509; A='Italia'; B='Italy'}
    @{Row=510; A='ahkera'; B='hardworking'}
    @{Row=511; A='koti-ikävä'; B='home sickness'}
    @{Row=512; A='banaani'; B='banana'}
    @{Row=513; A='olette'; B='(you plural) are'}
    @{Row=514; A='suolainen'; B='salty, savory'}
    @{Row=515; A='lattia'; B='floor (the type you stand on)'}
    @{Row=516; A='juotava'; B='drinkable'}
    @{Row=517; A='lehmä'; B='cow'}
    @{Row=518; A='mukava'; B='nice'}
    @{Row=519; A='istua'; B='to sit'}
    @{Row=520; A='virolainen'; B='Estonian person'}
    @{Row=521; A='melkein'; B='almost'}
    @{Row=522; A='sillä'; B='(it) has'}
    @{Row=523; A='juoda'; B='to drink'}
    @{Row=524; A='sana'; B='word'}
    @{Row=525; A='harmaa'; B='grey'}
    @{Row=526; A='väärä'; B='wrong (adjective)'}
    @{Row=527; A='seisoa'; B='to stand'}
    @{Row=528; A='pieni'; B='small'}
    @{Row=529; A='salaatti'; B='salad'}
    @{Row=530; A='kahvi'; B='coffee'}
    @{Row=531; A='paljonko'; B='what (time)'}
    @{Row=532; A='ihailla'; B='to admire'}
    @{Row=533; A='kerma'; B='cream'}
    @{Row=534; A='siivota'; B='to tidy up'}
    @{Row=535; A='keskellä'; B='in the middle of'}
    @{Row=536; A='missä'; B='where (relative)'}
    @{Row=537; A='nainen'; B='woman'}
    @{Row=538; A='no'; B='well (phrase)'}
    @{Row=539; A='ainakin'; B='at least'}
    @{Row=540; A='täydellinen'; B='perfect'}
    @{Row=541; A='puisto'; B='park'}
    @{Row=542; A='seitsemän'; B='seven'}
    @{Row=543; A='rouva'; B='Ms.'}
    @{Row=544; A='suo'; B='bog'}
    @{Row=545; A='pihvi'; B='steak'}
    @{Row=546; A='poni'; B='pony'}
    @{Row=547; A='ketsuppi'; B='ketchup'}
    @{Row=548; A='laulu'; B='song'}
    @{Row=549; A='seisoa'; B='to have stopped (clock, watch)'}
    @{Row=550; A='muu'; B='moo (cow sound)'}
    @{Row=551; A='tuoli'; B='chair'}
    @{Row=552; A='ratsastaa'; B='to ride (an animal)'}
    @{Row=553; A='tiskata'; B='to do the dishes'}
    @{Row=554; A='surullinen'; B='sad'}
    @{Row=555; A='sanoa'; B='to say (some sound)'}
    @{Row=556; A='rauha'; B='peace'}
    @{Row=557; A='munkki'; B='jelly doughnut'}
    @{Row=558; A='jäätelö'; B='ice cream'}
    @{Row=559; A='sinappi'; B='mustard'}
    @{Row=560; A='teillä'; B='you (plural) have'}
    @{Row=561; A='heillä'; B='they have'}
    @{Row=562; A='nätti'; B='pretty (adjective)'}
    @{Row=563; A='sulhanen'; B='groom'}
    @{Row=564; A='viineri'; B='Danish pastry'}
    @{Row=565; A='insinööri'; B='engineer'}
    @{Row=566; A='liha'; B='meat'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.A
    $ws.Cells.Item($u.Row, 2).Value = $u.B
}

$ws.Cells.Item(565, 3).Value = 1
$ws.Cells.Item(566, 3).Value = 1

